$d = $word.ActiveDocument

# Insert a new list-style paragraph after the current last paragraph,
# inheriting its paragraph/run formatting (ListParagraph style, numbering,
# sz/szCs run properties), then set its text.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs.Last
$newPara1.Range.Text = "Download icon on songs container"

# Repeat for the second new bullet item.
$newPara1.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Last
$newPara2.Range.Text = "In manipulateAPlaylist, addSongsFromSongBank, with an empty playlist their will be no image for it, however, when you add songs to it the image should appear. I need to make that happen."
